# add sheet 2 on fetch data
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Add the new worksheet right after Sheet1
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# Header row values for Sheet2 (row 2)
$headers = @("No", "Nama", "Token", "Dark", "Earth", "Fire", "Holy", "Ice", "Leaf", "Lightning", "Magma", "Metal", "Norm", "Rock", "Spirit", "Toxin", "Water", "Wind", "Wood")

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws2.Cells.Item(2, $i + 1).Value = $headers[$i]
}

# Reuse Sheet1's header-row formatting (bold, centered, filled, bordered) for the new header row
$ws1.Range("A2").Copy()
$destRange = $ws2.Range($ws2.Cells.Item(2, 1), $ws2.Cells.Item(2, $headers.Length))
$destRange.PasteSpecial(-4122)

# Column widths for Sheet2 (B..S), matching the source-file layout
$ws2.Columns.Item(2).ColumnWidth = 21.92
$ws2.Columns.Item(3).ColumnWidth = 19.25
$ws2.Columns.Item(4).ColumnWidth = 18.59
$ws2.Columns.Item(5).ColumnWidth = 18.42
$ws2.Columns.Item(6).ColumnWidth = 17.25
$ws2.Columns.Item(7).ColumnWidth = 15.75
$ws2.Columns.Item(8).ColumnWidth = 16.59
$ws2.Columns.Item(9).ColumnWidth = 17.42
$ws2.Columns.Item(10).ColumnWidth = 15.25
$ws2.Columns.Item(11).ColumnWidth = 14.92
$ws2.Columns.Item(12).ColumnWidth = 16.09
$ws2.Columns.Item(13).ColumnWidth = 16.25
$ws2.Columns.Item(14).ColumnWidth = 15.92
$ws2.Columns.Item(15).ColumnWidth = 16.09
$ws2.Columns.Item(16).ColumnWidth = 16.25
$ws2.Columns.Item(17).ColumnWidth = 16.09
$ws2.Columns.Item(18).ColumnWidth = 16.09
$ws2.Columns.Item(19).ColumnWidth = 16.09

# Selection / active sheet: leave Sheet1 active with A2 selected (as in the target file)
[void]$ws2.Range("A2").Select()
[void]$ws1.Activate()
[void]$ws1.Range("A2").Select()
